$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-TestRow($RowNum, $IdValue, $Height, $Title, $PreCondition, $Steps, $Expected) {

    $ws.Range("A$RowNum").Value = $IdValue

    $bCell = $ws.Range("B$RowNum")
    $bCell.Value = $Title
    $bCell.WrapText = $true
    $bCell.Font.Color = 0

    $cCell = $ws.Range("C$RowNum")
    $cCell.Value = $PreCondition
    $cCell.Font.Color = 0

    $dCell = $ws.Range("D$RowNum")
    $dCell.Value = $Steps
    $dCell.WrapText = $true
    $dCell.Font.Color = 0

    $eCell = $ws.Range("E$RowNum")
    $eCell.Value = $Expected
    $eCell.WrapText = $true
    $eCell.Font.Color = 0

    $ws.Rows.Item($RowNum).RowHeight = $Height
}

Add-TestRow 39 38 119 "Verify pagination on Home page" "Access home page" "1.) Load URL  
2.) Check for pagination elements at the bottom of home page
3.) Click on next and previous buttons" "Verify that all the page should be displayed and only 10 record are shown per page, and when next button is clicked subsequent 10 more records should be shown in accending order. When previousl button is access previous 10 record should be accessed"

Add-TestRow 40 39 68 "Verify end of  on Home page" "Access home page" "1.) Load URL  
2.) Check for pagination elements at the bottom of home page
3.) Click on next untill end of home page is reached " "Next button should be disabled"

Add-TestRow 41 40 68 "Verify beginning of  on Home page" "Access home page" "1.) Load URL  
2.) Check for pagination elements at the bottom of home page" "Previous button should be disabled when there are no more than 10 record and if the home page is at the bening of all records. "

Add-TestRow 42 41 51 " Verify pagenation on search value " "Access home page" "1.) Load URL  
2.) Search for a computer name which have less than 10 records" "Next and Previous button should not be displayed"

Add-TestRow 43 42 51 " Verify pagenation on search value with more than 10 records" "Access home page" "1.) Load URL  
2.) Search for a computer name which have more than 10 records" "Next  button should not be enabled, and by clicking it should display next set of records. "

# Update the sheet view: scroll so row 36 is at the top, and select A40:A43
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 36
$ws.Range("A40:A43").Select()
